# Update gh-pages to output generated at 456a3b4
# F2 (想去人数 / "want to go" count) of the first event row changes 1592 -> 1591
# F4 (想去人数) of the third event row changes 52 -> 54
# These updates apply identically to both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1591
    $ws.Range("F4").Value = 54
}
